$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sets": extend the time-index column (A) from 2 down to 24, i.e.
# rows 4..25 get values 3..24 (rows 2/3 already hold 1/2).
# ---------------------------------------------------------------------------
$sets = $wb.Worksheets.Item("Sets")
for ($i = 3; $i -le 24; $i++) {
    $sets.Cells.Item($i + 1, 1).Value = $i
}
$sets.Activate()
$sets.Range("A2:A25").Select()

# ---------------------------------------------------------------------------
# Sheet "General Data": bump the investment budget / capacity value (B3)
# from 100 to 2000.
# ---------------------------------------------------------------------------
$general = $wb.Worksheets.Item("General Data")
$general.Range("B3").Value = 2000
$general.Range("F10").Select()

# ---------------------------------------------------------------------------
# Sheet "Cost": no data changes, just move the remembered selection.
# ---------------------------------------------------------------------------
$cost = $wb.Worksheets.Item("Cost")
$cost.Range("C9").Select()

# ---------------------------------------------------------------------------
# Sheet "Demand": add a demand constraint value in every timestep (1..24).
# Timesteps 16..19 (rows 17..20) carry a demand of 22, every other
# timestep is 0.
# ---------------------------------------------------------------------------
$demand = $wb.Worksheets.Item("Demand")
$demandValues = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,22,22,22,22,0,0,0,0,0)
for ($i = 0; $i -lt $demandValues.Length; $i++) {
    $row = $i + 2
    $demand.Cells.Item($row, 1).Value = $i + 1
    $demand.Cells.Item($row, 2).Value = $demandValues[$i]
}
$demand.Activate()
$demand.Range("B17").Select()

# ---------------------------------------------------------------------------
# Sheet "irradiation": replace the flat 0.002 series with a per-timestep
# irradiation profile (B/C columns), then total column C in row 26.
# ---------------------------------------------------------------------------
$irr = $wb.Worksheets.Item("irradiation")
$irrValues = @(0,0,0,0,0,0,0,0.003,0.014,0.003,0.003,0,0.001,0.015,0,0,0,0,0,0,0,0,0,0)
for ($i = 0; $i -lt $irrValues.Length; $i++) {
    $row = $i + 2
    $irr.Cells.Item($row, 1).Value = $i + 1
    $irr.Cells.Item($row, 2).Value = $irrValues[$i]
    $irr.Cells.Item($row, 3).Value = $irrValues[$i]
}
$irr.Range("C26").Formula = "=SUM(C2:C25)"
$irr.Activate()
$irr.Range("C26").Select()

# ---------------------------------------------------------------------------
# Make "Demand" the active tab, matching the workbook's remembered view.
# ---------------------------------------------------------------------------
$demand.Activate()
$demand.Range("B17").Select()
